# Update localization status report for "Ready for handoff" generation.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns for zh-cn (E) and de-de (F), plus the
# "Latest HO Xliff Generate Date" column (G).
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-21 13:02:46"

# zh-cn sheet: Status (C) and Latest Handoff Datetime (H).
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-21 13:02:42"

# de-de sheet: Status (C) and Latest Handoff Datetime (H).
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-21 13:02:46"

# Widen the status columns slightly to fit the new "Ready for handoff" text.
$overview.Range("E1:F1").ColumnWidth = 16.33
$zhcn.Range("C1").ColumnWidth = 16.33
$dede.Range("C1").ColumnWidth = 16.33
